# SourceDictionary.xlsx update: align template with catalogue data model >4.x
# - Remove the "Repeated variables" sheet (its repeat-related columns are merged into "Variables")
# - Rewrite the header row of "Datasets", "Variables" and "Variable values" to the new column set

$wb = $excel.ActiveWorkbook

# --- Remove "Repeated variables" sheet -------------------------------------------------
$repeated = $wb.Worksheets.Item("Repeated variables")
$repeated.Delete() | Out-Null

# --- Datasets ----------------------------------------------------------------------------
$dsHeaders = @("resource","name","label","dataset type","unit of observation","keywords","description","number of rows","since version","until version")
$ws = $wb.Worksheets.Item("Datasets")
for ($i = 0; $i -lt $dsHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $dsHeaders[$i]
}
$ws.Range("A1:J1").EntireColumn.AutoFit() | Out-Null

# --- Variables -----------------------------------------------------------------------------
$varHeaders = @("resource","dataset","name","label","description","collection event","format","unit","since version","until version","repeat unit","repeat min","repeat max","example values","keywords","vocabularies","notes","useExternalDefinition.resource","useExternalDefinition.dataset","useExternalDefinition.name")
$ws = $wb.Worksheets.Item("Variables")
for ($i = 0; $i -lt $varHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $varHeaders[$i]
}
$ws.Range("A1:T1").EntireColumn.AutoFit() | Out-Null

# --- Variable values -------------------------------------------------------------------------
$valHeaders = @("resource","dataset","name","value","label","order","is missing","ontology term URI","since version","until version")
$ws = $wb.Worksheets.Item("Variable values")
for ($i = 0; $i -lt $valHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $valHeaders[$i]
}
$ws.Range("A1:J1").EntireColumn.AutoFit() | Out-Null

$wb.Worksheets.Item("Datasets").Activate() | Out-Null
